# Auto-generated edit script: update TPM-derived LR-pair statistics
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 31.824752
$ws.Range("H2").Value = 95.47425600000001
$ws.Range("I2").Value = 0.886907633630525
$ws.Range("J2").Value = 0.886907633630525
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 10.23061133333333
$ws.Range("N2").Value = 30.691834
$ws.Range("O2").Value = 0.4855635428718841
$ws.Range("P2").Value = 0.4855635428718841
$ws.Range("Q2").Value = 325.5866684917227
$ws.Range("R2").Value = 2930.280016425505
$ws.Range("S2").Value = 0.4306500127857567
$ws.Range("T2").Value = 0.4306500127857567
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 31.824752
$ws.Range("H3").Value = 95.47425600000001
$ws.Range("I3").Value = 0.886907633630525
$ws.Range("J3").Value = 0.886907633630525
$ws.Range("O3").Value = 0.4164864079521221
$ws.Range("P3").Value = 0.4164864079521222
$ws.Range("Q3").Value = 279.2681288121227
$ws.Range("R3").Value = 2513.413159309104
$ws.Range("S3").Value = 0.3693849745160941
$ws.Range("T3").Value = 0.3693849745160941
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 31.824752
$ws.Range("H4").Value = 95.47425600000001
$ws.Range("I4").Value = 0.886907633630525
$ws.Range("J4").Value = 0.886907633630525
$ws.Range("M4").Value = 2.034752
$ws.Range("N4").Value = 6.104255999999999
$ws.Range("O4").Value = 0.09657305490303886
$ws.Range("P4").Value = 0.09657305490303887
$ws.Range("Q4").Value = 64.755477781504
$ws.Range("R4").Value = 582.799300033536
$ws.Range("S4").Value = 0.08565137959652497
$ws.Range("T4").Value = 0.08565137959652498
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 31.824752
$ws.Range("H5").Value = 95.47425600000001
$ws.Range("I5").Value = 0.886907633630525
$ws.Range("J5").Value = 0.886907633630525
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.02901266666666667
$ws.Range("N5").Value = 0.087038
$ws.Range("O5").Value = 0.001376994272954919
$ws.Range("P5").Value = 0.001376994272954919
$ws.Range("Q5").Value = 0.9233209215253335
$ws.Range("R5").Value = 8.309888293728001
$ws.Range("S5").Value = 0.001221266732149232
$ws.Range("T5").Value = 0.001221266732149232
$ws.Range("I6").Value = 0.06502043684278042
$ws.Range("J6").Value = 0.06502043684278042
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 10.23061133333333
$ws.Range("N6").Value = 30.691834
$ws.Range("O6").Value = 0.4855635428718841
$ws.Range("P6").Value = 0.4855635428718841
$ws.Range("Q6").Value = 23.86921322219267
$ws.Range("R6").Value = 214.822918999734
$ws.Range("S6").Value = 0.03157155367245804
$ws.Range("T6").Value = 0.03157155367245804
$ws.Range("I7").Value = 0.06502043684278042
$ws.Range("J7").Value = 0.06502043684278042
$ws.Range("O7").Value = 0.4164864079521221
$ws.Range("P7").Value = 0.4164864079521222
$ws.Range("S7").Value = 0.02708012818412743
$ws.Range("T7").Value = 0.02708012818412744
$ws.Range("I8").Value = 0.06502043684278042
$ws.Range("J8").Value = 0.06502043684278042
$ws.Range("M8").Value = 2.034752
$ws.Range("N8").Value = 6.104255999999999
$ws.Range("O8").Value = 0.09657305490303886
$ws.Range("P8").Value = 0.09657305490303887
$ws.Range("Q8").Value = 4.747314481984
$ws.Range("R8").Value = 42.725830337856
$ws.Range("S8").Value = 0.006279222217037404
$ws.Range("T8").Value = 0.006279222217037405
$ws.Range("I9").Value = 0.06502043684278042
$ws.Range("J9").Value = 0.06502043684278042
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.02901266666666667
$ws.Range("N9").Value = 0.087038
$ws.Range("O9").Value = 0.001376994272954919
$ws.Range("P9").Value = 0.001376994272954919
$ws.Range("Q9").Value = 0.06768994581533334
$ws.Range("R9").Value = 0.6092095123380001
$ws.Range("S9").Value = 0.00008953276915753562
$ws.Range("T9").Value = 0.00008953276915753562
$ws.Range("G10").Value = 1.696588
$ws.Range("H10").Value = 5.089764000000001
$ws.Range("I10").Value = 0.04728133775640876
$ws.Range("J10").Value = 0.04728133775640876
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 10.23061133333333
$ws.Range("N10").Value = 30.691834
$ws.Range("O10").Value = 0.4855635428718841
$ws.Range("P10").Value = 0.4855635428718841
$ws.Range("Q10").Value = 17.35713242079733
$ws.Range("R10").Value = 156.214191787176
$ws.Range("S10").Value = 0.02295809387272402
$ws.Range("T10").Value = 0.02295809387272402
$ws.Range("G11").Value = 1.696588
$ws.Range("H11").Value = 5.089764000000001
$ws.Range("I11").Value = 0.04728133775640876
$ws.Range("J11").Value = 0.04728133775640876
$ws.Range("O11").Value = 0.4164864079521221
$ws.Range("P11").Value = 0.4164864079521222
$ws.Range("Q11").Value = 14.88787583089733
$ws.Range("R11").Value = 133.990882478076
$ws.Range("S11").Value = 0.01969203452533773
$ws.Range("T11").Value = 0.01969203452533774
$ws.Range("G12").Value = 1.696588
$ws.Range("H12").Value = 5.089764000000001
$ws.Range("I12").Value = 0.04728133775640876
$ws.Range("J12").Value = 0.04728133775640876
$ws.Range("M12").Value = 2.034752
$ws.Range("N12").Value = 6.104255999999999
$ws.Range("O12").Value = 0.09657305490303886
$ws.Range("P12").Value = 0.09657305490303887
$ws.Range("Q12").Value = 3.452135826176
$ws.Range("R12").Value = 31.069222435584
$ws.Range("S12").Value = 0.004566103227038787
$ws.Range("T12").Value = 0.004566103227038788
$ws.Range("G13").Value = 1.696588
$ws.Range("H13").Value = 5.089764000000001
$ws.Range("I13").Value = 0.04728133775640876
$ws.Range("J13").Value = 0.04728133775640876
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 0.6666666666666666
$ws.Range("M13").Value = 0.02901266666666667
$ws.Range("N13").Value = 0.087038
$ws.Range("O13").Value = 0.001376994272954919
$ws.Range("P13").Value = 0.001376994272954919
$ws.Range("Q13").Value = 0.04922254211466668
$ws.Range("R13").Value = 0.4430028790320001
$ws.Range("S13").Value = 0.00006510613130822202
$ws.Range("T13").Value = 0.00006510613130822202
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 0.3333333333333333
$ws.Range("G14").Value = 0.02836866666666667
$ws.Range("H14").Value = 0.085106
$ws.Range("I14").Value = 0.0007905917702857979
$ws.Range("J14").Value = 0.0007905917702857978
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 10.23061133333333
$ws.Range("N14").Value = 30.691834
$ws.Range("O14").Value = 0.4855635428718841
$ws.Range("P14").Value = 0.4855635428718841
$ws.Range("Q14").Value = 0.2902288027115556
$ws.Range("R14").Value = 2.612059224404
$ws.Range("S14").Value = 0.0003838825409453268
$ws.Range("T14").Value = 0.0003838825409453267
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 0.3333333333333333
$ws.Range("G15").Value = 0.02836866666666667
$ws.Range("H15").Value = 0.085106
$ws.Range("I15").Value = 0.0007905917702857979
$ws.Range("J15").Value = 0.0007905917702857978
$ws.Range("O15").Value = 0.4164864079521221
$ws.Range("P15").Value = 0.4164864079521222
$ws.Range("Q15").Value = 0.2489403360282222
$ws.Range("R15").Value = 2.240463024254
$ws.Range("S15").Value = 0.0003292707265628412
$ws.Range("T15").Value = 0.0003292707265628412
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 0.3333333333333333
$ws.Range("G16").Value = 0.02836866666666667
$ws.Range("H16").Value = 0.085106
$ws.Range("I16").Value = 0.0007905917702857979
$ws.Range("J16").Value = 0.0007905917702857978
$ws.Range("M16").Value = 2.034752
$ws.Range("N16").Value = 6.104255999999999
$ws.Range("O16").Value = 0.09657305490303886
$ws.Range("P16").Value = 0.09657305490303887
$ws.Range("Q16").Value = 0.05772320123733332
$ws.Range("R16").Value = 0.519508811136
$ws.Range("S16").Value = 0.00007634986243770104
$ws.Range("T16").Value = 0.00007634986243770104
$ws.Range("E17").Value = 1
$ws.Range("F17").Value = 0.3333333333333333
$ws.Range("G17").Value = 0.02836866666666667
$ws.Range("H17").Value = 0.085106
$ws.Range("I17").Value = 0.0007905917702857979
$ws.Range("J17").Value = 0.0007905917702857978
$ws.Range("K17").Value = 2
$ws.Range("L17").Value = 0.6666666666666666
$ws.Range("M17").Value = 0.02901266666666667
$ws.Range("N17").Value = 0.087038
$ws.Range("O17").Value = 0.001376994272954919
$ws.Range("P17").Value = 0.001376994272954919
$ws.Range("Q17").Value = 0.0008230506697777778
$ws.Range("R17").Value = 0.007407456028000001
$ws.Range("S17").Value = 0.000001088640339928834
$ws.Range("T17").Value = 0.000001088640339928834
